# Allowing for variable MCQ-s options.
#
# The "Q" header in C1 is renamed to "question", and the placeholder
# template rows (H4:H17) that used to show the literal text
# "The answer" now show "Option B" instead, so the answer key lines up
# with one of the actual MCQ options (matching rows 2/3 which already
# reference a real option letter instead of a fixed placeholder).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Q" column header to "question".
$ws.Range("C1").Value = "question"

# Point the template answer-key rows at "Option B" instead of the old
# fixed "The answer" placeholder.
$ws.Range("H4:H17").Value = "Option B"
